$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Ram
$ws.Range("B2").Value = "Ram"
$ws.Range("C2").Value = "Ram"
$ws.Range("E2").Formula = "=TRUE()"

# Row 3 - Gaurav
$ws.Range("B3").Value = "Gaurav"
$ws.Range("C3").Value = "Gaurav"
$ws.Range("D3").Value = 27
$ws.Range("E3").Formula = "=TRUE()"

# Row 4 - Anuj
$ws.Range("B4").Value = "Anuj"
$ws.Range("C4").Value = "Anuj"
$ws.Range("D4").Value = 28
$ws.Range("E4").Formula = "=TRUE()"

# Row 5 - Monu
$ws.Range("B5").Value = "Monu"
$ws.Range("C5").Value = "Monu"
$ws.Range("D5").Value = 29
$ws.Range("E5").Formula = "=FALSE()"

# Update selection to match the authored edit
$ws.Range("D2:D5").Select()
